# Converts a cell that currently holds a number into a text cell holding
# the given string, without leaving a residual style applied to the cell
# (Excel would otherwise keep a "@" text-format style on it).
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Overall": A2 (count of 990 filers) becomes text "2,051"
# ------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Range("A2") "2,051"

# ------------------------------------------------------------------
# Sheet "County": column B (No. of 990 Filers w/ Gov Grants) for every
# county row becomes text, and a new "Total" row (92) is appended.
# ------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
$countyCounts = @{
    2 = "2";   3 = "122"; 4 = "40";  5 = "1";   6 = "19";  7 = "8";
    8 = "2";   9 = "7";   10 = "21"; 11 = "5";  12 = "8";  13 = "1";
    14 = "10"; 15 = "13"; 16 = "25"; 17 = "10"; 18 = "34"; 19 = "14";
    20 = "62"; 21 = "4";  22 = "22"; 23 = "1";  24 = "5";  25 = "7";
    26 = "4";  27 = "17"; 28 = "4";  29 = "91"; 30 = "11"; 31 = "14";
    32 = "26"; 33 = "3";  34 = "18"; 35 = "12"; 36 = "14"; 37 = "7";
    38 = "3";  39 = "11"; 40 = "4";  41 = "24"; 42 = "14"; 43 = "34";
    44 = "12"; 45 = "19"; 46 = "108"; 47 = "15"; 48 = "19"; 49 = "488";
    50 = "11"; 51 = "6";  52 = "62"; 53 = "10"; 54 = "12"; 55 = "2";
    56 = "17"; 57 = "5";  58 = "5";  59 = "2";  60 = "3";  61 = "4";
    62 = "1";  63 = "46"; 64 = "6";  65 = "2";  66 = "5";  67 = "5";
    68 = "13"; 69 = "5";  70 = "6";  71 = "10"; 72 = "3";  73 = "108";
    74 = "5";  75 = "11"; 76 = "3";  77 = "10"; 78 = "52"; 79 = "3";
    80 = "4";  81 = "84"; 82 = "2";  83 = "41"; 84 = "15"; 85 = "1";
    86 = "9";  87 = "1";  88 = "23"; 89 = "9";  90 = "6";  91 = "8";
}
foreach ($r in $countyCounts.Keys) {
    Set-TextValue $wsCounty.Range("B$r") $countyCounts[$r]
}

# Append the new Total row at the bottom of the County sheet.
Set-TextValue $wsCounty.Range("A92") "Total"
Set-TextValue $wsCounty.Range("B92") "2,051"
Set-TextValue $wsCounty.Range("C92") "$3,305,071,391"
Set-TextValue $wsCounty.Range("D92") "10.03%"
Set-TextValue $wsCounty.Range("E92") "-11.39%"
Set-TextValue $wsCounty.Range("F92") "65.24%"

# ------------------------------------------------------------------
# Sheet "Congressional District": column B becomes text for every row,
# including the existing Total row (11), which uses "2,051".
# ------------------------------------------------------------------
$wsDistrict = $wb.Worksheets.Item("Congressional District")
$districtCounts = @{
    2 = "172"; 3 = "248"; 4 = "218"; 5 = "156"; 6 = "182";
    7 = "144"; 8 = "468"; 9 = "216"; 10 = "247";
}
foreach ($r in $districtCounts.Keys) {
    Set-TextValue $wsDistrict.Range("B$r") $districtCounts[$r]
}
Set-TextValue $wsDistrict.Range("B11") "2,051"

# ------------------------------------------------------------------
# Sheet "Size": column B becomes text for every row, including the
# existing Total row (8), which uses "2,051".
# ------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
$sizeCounts = @{
    2 = "677"; 3 = "534"; 4 = "351"; 5 = "141"; 6 = "275"; 7 = "73";
}
foreach ($r in $sizeCounts.Keys) {
    Set-TextValue $wsSize.Range("B$r") $sizeCounts[$r]
}
Set-TextValue $wsSize.Range("B8") "2,051"

# ------------------------------------------------------------------
# Sheet "Subsector": column B becomes text for every row, including
# the existing Total row (13), which uses "2,051".
# ------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")
$subsectorCounts = @{
    2 = "147"; 3 = "221"; 4 = "58"; 5 = "165"; 6 = "37"; 7 = "666";
    8 = "18";  9 = "200"; 10 = "55"; 11 = "458"; 12 = "26";
}
foreach ($r in $subsectorCounts.Keys) {
    Set-TextValue $wsSubsector.Range("B$r") $subsectorCounts[$r]
}
Set-TextValue $wsSubsector.Range("B13") "2,051"

Write-Host "Done applying text-formatting edits."
